# Applies the edits described by the commit diff:
#  1. Re-stamp the cached "datetimeFigureOut" date field text, on the
#     slide master and on every slide layout, from 1/30/2021 -> 2/4/2021.
#  2. Rework the "Amazon Elastic / Container Kubernetes" textbox on
#     slide 1 into a single "Amazon EKS" line (the shape uses
#     autofit-to-text, so its height shrinks on its own to match).

$p = $ppt.ActivePresentation

function Update-DatePlaceholders($shapes, $newDateText) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Type -eq 14 -and $shp.PlaceholderFormat.Type -eq 16) {
            $shp.TextFrame.TextRange.Text = $newDateText
        }
    }
}

# --- 1. Slide master + every custom (slide) layout ---
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes "2/4/2021"

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholders $layout.Shapes "2/4/2021"
}

# --- 2. "Amazon Elastic Container Kubernetes" -> "Amazon EKS" ---
$slide = $p.Slides.Item(1)
for ($k = 1; $k -le $slide.Shapes.Count; $k++) {
    $shp = $slide.Shapes.Item($k)
    if ($shp.Name -eq "TextBox 9") {
        $shp.TextFrame.TextRange.Text = "Amazon EKS"
    }
}
